$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.779.65'
$ws.Range('E2').Value = '  +0.83%  '
$ws.Range('D3').Value = '3.142.99'
$ws.Range('E3').Value = '  +1.11%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '586.69'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.32%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.54'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.36%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '3.138.34'
$ws.Range('E8').Value = '  +1.22%  '
$ws.Range('E9').Value = '  -0.37%  '
$ws.Range('E10').Value = '  +5.66%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.74'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.92%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.457'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.12%  '
$ws.Range('E13').Value = '  -0.48%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '37.04'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.50%  '
$ws.Range('D15').Value = '3.664.26'
$ws.Range('E15').Value = '  +1.16%  '
$ws.Range('E16').Value = '  -1.63%  '
$ws.Range('D17').Value = '3.143.47'
$ws.Range('E17').Value = '  +1.17%  '
$ws.Range('D18').Value = '63.582.47'
$ws.Range('E18').Value = '  +0.68%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.07'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.09%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '464.13'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.51%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.26'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.89%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.730'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.62%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.43'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.35%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.96'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.36%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '81.26'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.72%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.21'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.24%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.14%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.21'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +7.82%  '
$ws.Range('E29').Value = '  +0.43%  '
$ws.Range('B30').Value = 'FirstDigitalUSD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.01%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.21'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.44%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.95'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.00'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.110'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.49%  '
$ws.Range('E35').Value = '  -2.22%  '
$ws.Range('E36').Value = '  -0.70%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.29'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.62%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.32'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.45%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.00'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.15%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '51.00'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.99%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '440.73'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.56%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.80'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.07%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0371'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.56%  '
$ws.Range('D44').Value = '2.905.63'
$ws.Range('E44').Value = '  -0.26%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.279'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.107'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.91%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '36.90'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.15%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '125.40'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.80%  '
$ws.Range('E50').Value = '  -0.96%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '24.28'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.17%  '
